$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Limit the number of cases to 100 (free tier): clear the "Cases" values
# in column A for rows 102 through 504, keeping the first 100 data rows
# (rows 2-101) intact.
$ws.Range("A102:A504").ClearContents()

# Update the selected cell to match the new state recorded in the workbook.
$ws.Range("B20").Select()
